$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98, shifting rows 98:187 down to 99:188
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with the new weekly record
$ws.Cells.Item(98, 1).Value2 = 3
$ws.Cells.Item(98, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(98, 3).Value2 = "Coquimbo"
$ws.Cells.Item(98, 4).Value2 = 44589
$ws.Cells.Item(98, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat
$ws.Cells.Item(98, 5).Value2 = 5
$ws.Cells.Item(98, 6).Value2 = "Fruta"
$ws.Cells.Item(98, 7).Value2 = 100101
$ws.Cells.Item(98, 8).Value2 = "Berries"
$ws.Cells.Item(98, 9).Value2 = 100101001
$ws.Cells.Item(98, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(98, 11).Value2 = "Sin especificar"
$ws.Cells.Item(98, 12).Value2 = "Primera"
$ws.Cells.Item(98, 13).Value2 = 310
$ws.Cells.Item(98, 14).Value2 = 4500
$ws.Cells.Item(98, 15).Value2 = 5000
$ws.Cells.Item(98, 16).Value2 = 4758
$ws.Cells.Item(98, 17).Value2 = "$/bandeja 2 kilos"
$ws.Cells.Item(98, 18).Value2 = "Provincia de Linares"
$ws.Cells.Item(98, 19).Value2 = 2379
$ws.Cells.Item(98, 20).Value2 = 2
